$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 ("Tours have computed attributes"): recorded points 1 -> 2
$ws.Range("D29").Value = 2

# Row 67 ("Describes app architecture..."): replace the "x" placeholder with the actual points value (3)
$ws.Range("D67").Value = 3

# Rows 70-75 (protocol items): replace the "x" placeholder with the actual points value (1 each)
$ws.Range("D70:D75").Value = 1

# Move/record the active selection from A18 to E18
$ws.Range("E18").Select()
